$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (VBA RGB(r,g,b) = r + g*256 + b*65536)
$colorGreen  = 65280    # 00FF00
$colorYellow = 65535    # FFFF00

$xlCenter = -4108

# --- Row 3: mode 0 -> 8, "Portcullis  - Not implemented" -> "Corner Shot" ---
$a3 = $ws.Range("A3")
$a3.Value2 = 8
$a3.Font.Bold = $true
$a3.Font.Color = $colorGreen
$a3.HorizontalAlignment = $xlCenter

$b3 = $ws.Range("B3")
$b3.Value2 = "Corner Shot"
$b3.Font.Bold = $true
$b3.Font.Color = $colorGreen

# --- Row 4: mode 1 -> 100, "Cheval De Frise " -> "Low Bar w gyro" ---
$a4 = $ws.Range("A4")
$a4.Value2 = 100
$a4.Font.Bold = $true
$a4.Font.Color = $colorGreen
$a4.HorizontalAlignment = $xlCenter

$b4 = $ws.Range("B4")
$b4.Value2 = "Low Bar w gyro"
$b4.Font.Bold = $true
$b4.Font.Color = $colorGreen

# --- Row 5: mode 2 -> 101, "Moat " -> "Portcullis  w gyro" (highlighted yellow) ---
$a5 = $ws.Range("A5")
$a5.Value2 = 101
$a5.Font.Bold = $true
$a5.Font.Color = $colorYellow
$a5.HorizontalAlignment = $xlCenter

$b5 = $ws.Range("B5")
$b5.Value2 = "Portcullis  w gyro"
$b5.Font.Bold = $true
$b5.Font.Color = $colorYellow

# --- Row 6: mode 3 -> 102, "Ramparts " -> "Cheval De Frise w gyro" ---
$a6 = $ws.Range("A6")
$a6.Value2 = 102
$a6.Font.Bold = $true
$a6.Font.Color = $colorGreen
$a6.HorizontalAlignment = $xlCenter

$b6 = $ws.Range("B6")
$b6.Value2 = "Cheval De Frise w gyro"
$b6.Font.Bold = $true
$b6.Font.Color = $colorGreen

# --- Row 7: mode 4 -> 103, "Drawbridge - Not implemented" -> "Moat w gyro" ---
# (keeps the original red-highlighted style, unlike the rows above)
$a7 = $ws.Range("A7")
$a7.Value2 = 103

$b7 = $ws.Range("B7")
$b7.Value2 = "Moat w gyro"

# --- Row 8: mode 5 -> 200, "Sally Port  - Not implemented" -> "Low Bar w gyro, 2 ball w spy bot" (yellow) ---
$a8 = $ws.Range("A8")
$a8.Value2 = 200
$a8.Font.Bold = $true
$a8.Font.Color = $colorYellow
$a8.HorizontalAlignment = $xlCenter

$b8 = $ws.Range("B8")
$b8.Value2 = "Low Bar w gyro, 2 ball w spy bot"
$b8.Font.Bold = $true
$b8.Font.Color = $colorYellow

# --- Row 9: mode 6 -> "default", "Rock Wall  - Not implemented" -> "Corner Shot" ---
$a9 = $ws.Range("A9")
$a9.Value2 = "default"
$a9.Font.Bold = $true
$a9.Font.Color = $colorGreen
$a9.HorizontalAlignment = $xlCenter

$b9 = $ws.Range("B9")
$b9.Value2 = "Corner Shot"
$b9.Font.Bold = $true
$b9.Font.Color = $colorGreen

# --- Remove old rows 10-14 entirely (content + formatting) ---
$ws.Range("A10:B14").Clear()

# --- Refresh the stale selection/used-range so it matches the new extent ---
[void]$ws.Range("A1:B9").Select()
